# "Signed Off time sheets"
#
# The supervisor (Ankita Gangotra) has signed off on Prakruti Sinha's
# weekly timesheet:
#   - Supervisor Name is filled in on the header (merged G6:I6)
#   - Supervisor Signature block is completed: initials in A27 (merged
#     A27:C27) and the signature date in D27 (merged D27:E27), formatted
#     like the neighbouring "Date" cell (D25) used for the employee
#     signature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name: ...
$ws.Range("G6").Value = "Ankita Gangotra"

# Supervisor Signature  |  Date
$ws.Range("A27").Value = "A.G"
$ws.Range("D27").NumberFormat = "mm-dd-yy"
$ws.Range("D27").Value = 41800

# Leave the selection where the sign-off happened.
$ws.Range("D27:E27").Select()
